$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.441.12'
$ws.Range("E2").Value = '  -3.10%  '
$ws.Range("D3").Value = '1.991.24'
$ws.Range("E3").Value = '  -4.88%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.016'
$ws.Range("D4").Style = $ws.Range("C2").Style
$ws.Range("E4").Value = '  +1.33%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '328.34'
$ws.Range("D5").Style = $ws.Range("C2").Style
$ws.Range("E5").Value = '  -4.22%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4983'
$ws.Range("D7").Style = $ws.Range("C2").Style
$ws.Range("E7").Value = '  -4.76%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4202'
$ws.Range("D8").Style = $ws.Range("C2").Style
$ws.Range("E8").Value = '  -4.88%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '53.77'
$ws.Range("D9").Style = $ws.Range("C2").Style
$ws.Range("E9").Value = '  -1.50%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08922'
$ws.Range("D10").Style = $ws.Range("C2").Style
$ws.Range("E10").Value = '  -4.14%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.109'
$ws.Range("D11").Style = $ws.Range("C2").Style
$ws.Range("E11").Value = '  -4.99%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '23.16'
$ws.Range("D12").Style = $ws.Range("C2").Style
$ws.Range("E12").Value = '  -6.66%  '
$ws.Range("D13").Value = '2.002.41'
$ws.Range("E13").Value = '  -2.98%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.944'
$ws.Range("D14").Style = $ws.Range("C2").Style
$ws.Range("E14").Value = '  -7.47%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.418'
$ws.Range("D15").Style = $ws.Range("C2").Style
$ws.Range("E15").Value = '  -7.07%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.016'
$ws.Range("D16").Style = $ws.Range("C2").Style
$ws.Range("E16").Value = '  +1.38%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '93.43'
$ws.Range("D17").Style = $ws.Range("C2").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001104'
$ws.Range("D18").Style = $ws.Range("C2").Style
$ws.Range("E18").Value = '  -4.68%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06667'
$ws.Range("D19").Style = $ws.Range("C2").Style
$ws.Range("E19").Value = '  -0.02%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.40'
$ws.Range("D20").Style = $ws.Range("C2").Style
$ws.Range("E20").Value = '  -8.33%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.014'
$ws.Range("D21").Style = $ws.Range("C2").Style
$ws.Range("E21").Value = '  +1.28%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.908'
$ws.Range("D22").Style = $ws.Range("C2").Style
$ws.Range("E22").Value = '  -6.77%  '
$ws.Range("D23").Value = '29.488.40'
$ws.Range("E23").Value = '  -2.99%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.92'
$ws.Range("D24").Style = $ws.Range("C2").Style
$ws.Range("E24").Value = '  -4.88%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.302'
$ws.Range("D25").Style = $ws.Range("C2").Style
$ws.Range("E25").Value = '  -0.49%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '157.40'
$ws.Range("D26").Style = $ws.Range("C2").Style
$ws.Range("E26").Value = '  -3.53%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.63'
$ws.Range("D27").Style = $ws.Range("C2").Style
$ws.Range("E27").Value = '  -5.45%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.237'
$ws.Range("D28").Style = $ws.Range("C2").Style
$ws.Range("E28").Value = '  -8.32%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.287'
$ws.Range("D29").Style = $ws.Range("C2").Style
$ws.Range("E29").Value = '  -8.79%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '127.21'
$ws.Range("D30").Style = $ws.Range("C2").Style
$ws.Range("E30").Value = '  -4.47%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.047'
$ws.Range("D31").Style = $ws.Range("C2").Style
$ws.Range("E31").Value = '  -7.83%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09879'
$ws.Range("D32").Style = $ws.Range("C2").Style
$ws.Range("E32").Value = '  -5.58%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.547'
$ws.Range("D33").Style = $ws.Range("C2").Style
$ws.Range("E33").Value = '  -6.41%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.810'
$ws.Range("D34").Style = $ws.Range("C2").Style
$ws.Range("E34").Value = '  -1.25%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.795'
$ws.Range("D35").Style = $ws.Range("C2").Style
$ws.Range("E35").Value = '  -7.35%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02449'
$ws.Range("D36").Style = $ws.Range("C2").Style
$ws.Range("E36").Value = '  -6.95%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '9.202'
$ws.Range("D37").Style = $ws.Range("C2").Style
$ws.Range("E37").Value = '  -9.16%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.295'
$ws.Range("D38").Style = $ws.Range("C2").Style
$ws.Range("E38").Value = '  -3.38%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06346'
$ws.Range("D39").Style = $ws.Range("C2").Style
$ws.Range("E39").Value = '  -7.19%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6517'
$ws.Range("D40").Style = $ws.Range("C2").Style
$ws.Range("E40").Value = '  -6.62%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '11.56'
$ws.Range("D41").Style = $ws.Range("C2").Style
$ws.Range("E41").Value = '  -8.05%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.2030'
$ws.Range("D42").Style = $ws.Range("C2").Style
$ws.Range("E42").Value = '  -8.29%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.014'
$ws.Range("D43").Style = $ws.Range("C2").Style
$ws.Range("E43").Value = '  +1.37%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6288'
$ws.Range("D44").Style = $ws.Range("C2").Style
$ws.Range("E44").Value = '  -7.59%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.33'
$ws.Range("D45").Style = $ws.Range("C2").Style
$ws.Range("E45").Value = '  -7.56%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.190'
$ws.Range("D46").Style = $ws.Range("C2").Style
$ws.Range("E46").Value = '  -6.62%  '
$ws.Range("E47").Value = '  -5.52%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.502'
$ws.Range("D48").Style = $ws.Range("C2").Style
$ws.Range("E48").Value = '  -3.59%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00000000335'
$ws.Range("D49").Style = $ws.Range("C2").Style
$ws.Range("E49").Value = '  -4.64%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06948'
$ws.Range("D50").Style = $ws.Range("C2").Style
$ws.Range("E50").Value = '  -4.08%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.120'
$ws.Range("D51").Style = $ws.Range("C2").Style
$ws.Range("E51").Value = '  -9.33%  '
